# Update the regex-pattern column (A) with refined/escaped patterns.
# The category column (B) is unchanged; only some A values were rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "(#NUM )؟(درصد |واحد )؟(ضرر|سود)"
$ws.Range("A3").Value2 = "(اطلاعیه (ی )؟|آگهی |اعلامیه )؟((افشا (ی )؟(اطلاعات )؟(الف|ب|با اهمیت))|ثبت افزایش سرمایه|صورت (ها ی )؟مالی|فعالیت ماهانه|دعوت به مجمع( عمومی| عادی)؟|پذیره نویسی عمومی)"
$ws.Range("A5").Value2 = "#NUM (درصد |واحد )؟افزایش سرمایه"
$ws.Range("A7").Value2 = "دامنه (ی )؟نوسان"
$ws.Range("A10").Value2 = "(#NUM )؟(درصد |واحد )؟(تاثیر )؟(مثبت|منفی)"

# Update the active cell selection to A3 (was B12)
$ws.Range("A3").Select()
